$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (rows 2-15) after the "Aggiorno file need_to_buy.xlsx da R" update.
# Columns: A=Giorno, B=TISG_PDR_G, C=fcs, D=buy_BEE_MWH, E=sell_lago_MWH, F=need_to_buy_MW
$data = @(
    @(45943, 5582.50657861397, 5419.07112850915, 3620, 8138.218348, 181.449287412299),
    @(45944, 5564.50698273871, 5674.47581712338, 3620, 8104.081383, 191.418759057695),
    @(45945, 5816.11417805593, 5930.05926388495, 3620, 8539.482492, 209.726149076209),
    @(45946, 6218.61815962666, 6240.58263879982, 3620, 9246.004617, 235.332045673881),
    @(45947, 6118.17196890358, 5582.19307855565, 3620, 9082.81838,  205.284978735503),
    @(45948, 1742.27770790123, 3952.20454120997, 3620, 4473.856137, 127.657623762864),
    @(45949, 1636.94065696827, 3964.93507613243, 3620, 4358.695036, 127.778727298507),
    @(45950, 6392.95297294923, 6359.019888416,   3620, 9928.525459, 261.441348936115),
    @(45951, 6392.95297294923, 6565.60795547408, 3620, 9928.525459, 270.049185063535),
    @(45952, 6392.95297294923, 6466.92505638208, 3620, 9928.525459, 265.937397601369),
    @(45953, 6392.95297294923, 6391.1291919267,  3620, 9928.525459, 262.779236582395),
    @(45954, 6392.95297294923, 5684.29390541131, 3620, 9928.525459, 233.32776631092),
    @(45955, 2044.14689306201, 3810.81046944629, 3620, 5157.785259, 137.685368141012),
    @(45956, 1933.9611505005,  3912.78100482224, 3620, 5037.33837,  141.506592680073)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
